$d = $word.ActiveDocument

# --- Step 1: remove stray _GoBack bookmark from the title line ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 2: collapse the ">>>  your stuff..." paragraph runs into one,
#     dropping the proofErr markers, via Find & Replace with identical text ---
$d.Content.Find.Execute(
    ">>>  your stuff after this line >>>", $false, $false, $false, $false, $false,
    $true, 1, $false, ">>>  your stuff after this line >>>", 2
) | Out-Null

# --- Step 3: append four new paragraphs after "Ben changing things up!" ---
$texts = @(
    "Modification of the file can easily done by simple form. Today I will tell something about myself rather commenting about other topics.",
    "Basically, being an IT student I do love to spend time cooking as well. I love driving as well. I mostly travel within the state and also sometimes outside driving my own car with my friends and family.",
    "I like to go to the Pub and drink and dance with friends.",
    "Cheers."
)

$idx = 5
foreach ($t in $texts) {
    $d.Paragraphs($idx).Range.InsertParagraphAfter()
    $idx = $idx + 1
    $d.Paragraphs($idx).Range.Text = $t
}

# --- Step 4: re-add the _GoBack bookmark right after "Cheers." (before its
#     paragraph mark). A zero-length range sitting exactly on a paragraph-end
#     boundary cannot receive a fresh bookmark, so temporarily pad with one
#     extra character, bookmark just before it, then remove the padding. ---
$cheers = $d.Paragraphs($idx).Range
$padPos = $cheers.End - 1
$d.Range($padPos, $padPos).InsertAfter("Z") | Out-Null

$bmPos = $padPos
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos)) | Out-Null

$d.Range($bmPos, $bmPos + 1).Delete() | Out-Null
